$d = $word.ActiveDocument

# Locate the three paragraphs to remove:
#   1) the blank paragraph right after "LOB1256: ..."
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "(c) 2020 . Contact: ..." paragraph
# They sit between the "LOB1256..." requisito paragraph and the trailing
# blank / page-break paragraph that precedes the section break.

$paras = $d.Paragraphs
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $paras.Count; $i++) {
    $text = $paras.Item($i).Range.Text
    if ($text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startIndex = $i - 1
    }
    if ($text -like "*Powered by Jekyll and Github pages*") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $rngStart = $paras.Item($startIndex).Range.Start
    $rngEnd = $paras.Item($endIndex).Range.End
    $rng = $d.Range($rngStart, $rngEnd)
    $rng.Delete()
}
